# Applies "added more attributes to particle system import export":
#  - row 133: now documents the "color gradient" (iaGradientColor4f) field
#  - row 134: now documents the "emission gradient" (iaGradientf) field
#  - rows 137-141: new vortex torque/range fields added
#  - rows 143-145: the Switch-Chunk / LOD-Chunk legend entries move down to 150-152
#  - rows 150-152: Switch-Chunk / LOD-Chunk legend entries (previously at 143-145)
#  - rows 217-223: sheet grows by 7 blank rows
#  - dimension / sheetView follow the sheet's new extent & selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Row 133: emitterID/node id of emitter -> color gradient/iaGradientColor4f ----
$ws.Range("B133").Value = "color gradient"

$ws.Range("C133").Value = "iaGradientColor4f"
$ws.Range("C133").WrapText = $false

$ws.Range("D133").Value = "gradient"
$ws.Range("D133").WrapText = $false

$ws.Range("E133").Value = "variable"
$ws.Range("E133").WrapText = $false

$ws.Range("F133").Value = "na"

# ---- Row 134: color gradient/iaGradientColor4f -> emission gradient/iaGradientf ----
$ws.Range("B134").Value = "emission gradient"
$ws.Range("C134").Value = "iaGradientf"
# D134/E134/F134 (gradient / variable / na) stay as they were

# ---- Rows 137-141: new vortex torque/range attributes ----
$ws.Range("B137").Value = "min vortex torque"
$ws.Range("D137").Value = "float"
$ws.Range("E137").Value = 4
$ws.Range("F137").Value = 0

$ws.Range("B138").Value = "max vortex torque"
$ws.Range("D138").Value = "float"
$ws.Range("E138").Value = 4
$ws.Range("F138").Value = 0

$ws.Range("B139").Value = "min vortex range"
$ws.Range("D139").Value = "float"
$ws.Range("E139").Value = 4

$ws.Range("B140").Value = "max vortex range"
$ws.Range("D140").Value = "float"
$ws.Range("E140").Value = 4

$ws.Range("B141").Value = "vortex check range"
$ws.Range("D141").Value = "uint"
$ws.Range("E141").Value = 1

# ---- Rows 143 & 145: clear the Switch-Chunk / LOD-Chunk legend (moved to 150/152) ----
$ws.Range("A143").ClearContents()
$ws.Range("C143").ClearContents()
$ws.Range("D143").ClearContents()
$ws.Cells.Item(143, 1).EntireRow.AutoFit()

$ws.Range("A145").ClearContents()
$ws.Range("C145").ClearContents()
$ws.Range("D145").ClearContents()
$ws.Cells.Item(145, 1).EntireRow.AutoFit()

# ---- Rows 150 & 152: write the Switch-Chunk / LOD-Chunk legend here instead ----
$ws.Range("A150").Value = "Switch Chunk"
$ws.Range("C150").Value = "* just contains children`n* gives hint to the application how to handle the children`n* there can only be one active child for rendering"
$ws.Range("D151").WrapText = $true
$ws.Rows.Item(150).RowHeight = 42.75

$ws.Range("A152").Value = "LOD-Chunk"
$ws.Range("C152").Value = "similar to switch chunk but the application can decide which LODs to show even in parallel"
$ws.Rows.Item(152).RowHeight = 28.5

# ---- Rows 217-223: sheet grows with 7 new blank rows (style like the existing tail rows) ----
for ($r = 217; $r -le 223; $r++) {
    $ws.Range("C$r").WrapText = $true
}

# ---- Sheet view / dimension bookkeeping ----
$ws.Range("A1:G223").Select()
$ws.Application.ActiveWindow.ScrollRow = 115
$ws.Range("C136").Select()
